$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (89b0dd29-...) status + datetime updated to "Ready for handoff" ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-10-18 11:45:31"

# --- zh-cn sheet: row 3 (89b0dd29-...) status + handoff datetime + error detail ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("H3").Value = "2016-10-18 11:45:20"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0ad011ea353a38514bb4ba5016fa94de9e555702/e2e/89b0dd29-cb79-48f1-9fdc-c0e62a8961ea.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c5591b558086a8402f017fefda6a0e94674694c/e2e/89b0dd29-cb79-48f1-9fdc-c0e62a8961ea.md."
$zh.Columns.Item(16).ColumnWidth = 39.1428571428571

# --- de-de sheet: row 3 (89b0dd29-...) status + handoff datetime + error detail ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("H3").Value = "2016-10-18 11:45:31"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0ad011ea353a38514bb4ba5016fa94de9e555702/e2e/89b0dd29-cb79-48f1-9fdc-c0e62a8961ea.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c5591b558086a8402f017fefda6a0e94674694c/e2e/89b0dd29-cb79-48f1-9fdc-c0e62a8961ea.md."
$de.Columns.Item(16).ColumnWidth = 39.1428571428571
